$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

function Rename-SimulationBlock($startRow, $endRow, $newName) {
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $newName
    }
}

# Cultivar/treatment names reworked (TT targets adjusted for Emerald cultivar).
# The write order below controls the first-appearance order of the new
# shared strings so it matches the authored workbook.
Rename-SimulationBlock 2   59  "exp4WaterIrrCultivarEmerald"
Rename-SimulationBlock 122 172 "exp4WaterRF_IrrCultivarEmerald"
Rename-SimulationBlock 173 218 "exp4SowSoybeanIrrCultivarEmerald"
Rename-SimulationBlock 272 310 "exp4SowSoybeanRFIrrCultivarEmerald"
Rename-SimulationBlock 60  121 "exp4WaterRFCultivarEmerald"
Rename-SimulationBlock 219 271 "exp4SowSoybeanRFCultivarEmerald"

# Update the saved view state: scroll so row 296 becomes the top row and the
# active selection is A316 (graphs for phenology/node numbers were being
# reviewed near the bottom of the sheet).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 296
$ws.Range("A316").Select()
